$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/15/2025  Through  9/21/2025"

# --- Crime-data table updates (rows 14-33) ---
# Row 14
$ws.Range("C23").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("C14").Value = "0"
$ws.Range("N14").Value = -86.95652173913

# Row 15
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 29
$ws.Range("K15").Value = 3.571428571428
$ws.Range("L15").Value = 20.833333333333
$ws.Range("M15").Value = 61.111111111111
$ws.Range("N15").Value = 26.086956521739

# Row 16
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 16
$ws.Range("E16").Value = -37.5
$ws.Range("F16").Value = 36
$ws.Range("G16").Value = 51
$ws.Range("H16").Value = -29.411764705882
$ws.Range("I16").Value = 296
$ws.Range("J16").Value = 385
$ws.Range("K16").Value = -23.116883116883
$ws.Range("L16").Value = -12.166172106824
$ws.Range("M16").Value = 14.728682170542
$ws.Range("N16").Value = -74.080560420315

# Row 17
$ws.Range("D17").Value = 13
$ws.Range("E17").Value = -15.384615384615
$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 57
$ws.Range("H17").Value = -24.561403508771
$ws.Range("I17").Value = 422
$ws.Range("J17").Value = 607
$ws.Range("K17").Value = -30.477759472817
$ws.Range("L17").Value = -22.568807339449
$ws.Range("M17").Value = 126.881720430108
$ws.Range("N17").Value = 17.548746518105

# Row 18
$ws.Range("C23").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = "0"
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 8
$ws.Range("H18").Value = -55.555555555555
$ws.Range("I18").Value = 121
$ws.Range("J18").Value = 185
$ws.Range("K18").Value = -34.594594594594
$ws.Range("L18").Value = -13.571428571428
$ws.Range("M18").Value = -42.654028436019
$ws.Range("N18").Value = -92.512376237623

# Row 19
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 27
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 70
$ws.Range("G19").Value = 91
$ws.Range("H19").Value = -23.076923076923
$ws.Range("I19").Value = 571
$ws.Range("J19").Value = 825
$ws.Range("K19").Value = -30.787878787878
$ws.Range("L19").Value = -26.512226512226
$ws.Range("M19").Value = 45.663265306122
$ws.Range("N19").Value = -35.260770975056

# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 26
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 180
$ws.Range("J20").Value = 219
$ws.Range("K20").Value = -17.808219178082
$ws.Range("L20").Value = -26.829268292682
$ws.Range("M20").Value = 63.636363636363
$ws.Range("N20").Value = -88.882025941939

# Row 21
$ws.Range("C21").Value = 45
$ws.Range("D21").Value = 64
$ws.Range("E21").Value = -29.6875
$ws.Range("F21").Value = 186
$ws.Range("G21").Value = 244
$ws.Range("H21").Value = -23.770491803278
$ws.Range("I21").Value = 1622
$ws.Range("J21").Value = 2251
$ws.Range("K21").Value = -27.943136383829
$ws.Range("L21").Value = -21.680347658136
$ws.Range("M21").Value = 37.807986406117
$ws.Range("N21").Value = -71.362994350282

# Row 22
$ws.Range("C23").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = "0"
$ws.Range("I29").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("K29").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 4
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 39
$ws.Range("J22").Value = 39
$ws.Range("L22").Value = 5.405405405405
$ws.Range("M22").Value = 62.5

# Row 24
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = 2.083333333333
$ws.Range("F24").Value = 204
$ws.Range("G24").Value = 209
$ws.Range("H24").Value = -2.392344497607
$ws.Range("I24").Value = 1755
$ws.Range("J24").Value = 2285
$ws.Range("K24").Value = -23.194748358862
$ws.Range("L24").Value = -14.598540145985
$ws.Range("M24").Value = 38.735177865612

# Row 25
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 32
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 120
$ws.Range("G25").Value = 153
$ws.Range("H25").Value = -21.56862745098
$ws.Range("I25").Value = 1245
$ws.Range("J25").Value = 1810
$ws.Range("K25").Value = -31.215469613259
$ws.Range("L25").Value = -16.944629753168

# Row 26
$ws.Range("C26").Value = 31
$ws.Range("D26").Value = 39
$ws.Range("E26").Value = -20.51282051282
$ws.Range("F26").Value = 108
$ws.Range("G26").Value = 117
$ws.Range("H26").Value = -7.692307692307
$ws.Range("I26").Value = 777
$ws.Range("J26").Value = 1008
$ws.Range("K26").Value = -22.916666666666
$ws.Range("L26").Value = -8.372641509433
$ws.Range("M26").Value = 67.096774193548

# Row 27
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 200
$ws.Range("L27").Value = 2.631578947368

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 6
$ws.Range("E28").Value = -66.666666666666
$ws.Range("G28").Value = 17
$ws.Range("H28").Value = -23.529411764705
$ws.Range("I28").Value = 87
$ws.Range("J28").Value = 112
$ws.Range("K28").Value = -22.321428571428
$ws.Range("L28").Value = -14.705882352941

# Row 29
$ws.Range("N29").Value = -85.416666666666

# Row 30
$ws.Range("N30").Value = -90.90909090909

# Row 31
$ws.Range("I29").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("D31").Value = 1
$ws.Range("K29").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100
$ws.Range("I29").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("G31").Value = 1
$ws.Range("K29").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 3
$ws.Range("K31").Value = 133.333333333333

# Row 33
$ws.Range("C23").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C33").Value = "0"
$ws.Range("C23").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("D33").Value = "0"
$ws.Range("E23").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("E33").Value = "***.*"
$ws.Range("F33").Value = 1
$ws.Range("H33").Value = 0

$excel.CutCopyMode = 0
